$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 13
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "rwef"
$ws.Range("D8").Value = "wefwef"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = "'"
$ws.Range("F8").Style = "Normal"
